# Rearrange the flux-box data table:
#  - Insert a new blank column before column B.
#  - Move the old last column (now shifted to Q) into the new blank column B
#    (this is the "Coarse_seds_subsurface" label + the per-row flag values).
#  - Delete the now-empty column Q left behind by the move.
#  - Update the active selection to reflect where the user ended up (E13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new empty column at B; everything B:P shifts right to C:Q.
$ws.Range("B1").EntireColumn.Insert() | Out-Null

# 2. Cut the data that used to be column P (now at Q after the shift) and
#    paste it into the freshly inserted column B.
$ws.Range("Q1:Q12").Cut($ws.Range("B1:B12")) | Out-Null

# 3. Remove the leftover empty column Q.
$ws.Range("Q1").EntireColumn.Delete() | Out-Null

# 4. Leave the selection where the editor ended up.
$ws.Range("E13").Select() | Out-Null
